$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text columns (Coin name / Link) -- plain text, no numeric ambiguity
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('B39').Value = 'EnergySwap'
$ws.Range('C39').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'

# Volume/percentage column (E) -- always textual due to '%' and spacing
$ws.Range('E2').Value = '  -0.39%  '
$ws.Range('E3').Value = '  +2.21%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('E5').Value = '  -0.35%  '
$ws.Range('E6').Value = '  +14.05%  '
$ws.Range('E7').Value = '  +0.69%  '
$ws.Range('E8').Value = '  +2.92%  '
$ws.Range('E9').Value = '  +1.86%  '
$ws.Range('E10').Value = '  -0.01%  '
$ws.Range('E11').Value = '  +2.28%  '
$ws.Range('E12').Value = '  +2.83%  '
$ws.Range('E13').Value = '  +0.28%  '
$ws.Range('E14').Value = '  +3.59%  '
$ws.Range('E15').Value = '  +2.22%  '
$ws.Range('E16').Value = '  +5.01%  '
$ws.Range('E17').Value = '  -0.48%  '
$ws.Range('E18').Value = '  +14.26%  '
$ws.Range('E19').Value = '  +2.39%  '
$ws.Range('E20').Value = '  +0.49%  '
$ws.Range('E21').Value = '  +1.20%  '
$ws.Range('E22').Value = '  -1.00%  '
$ws.Range('E23').Value = '  +2.47%  '
$ws.Range('E25').Value = '  +1.36%  '
$ws.Range('E26').Value = '  -0.93%  '
$ws.Range('E27').Value = '  +5.43%  '
$ws.Range('E28').Value = '  -1.76%  '
$ws.Range('E29').Value = '  +8.22%  '
$ws.Range('E30').Value = '  +6.41%  '
$ws.Range('E31').Value = '  -0.73%  '
$ws.Range('E32').Value = '  +0.13%  '
$ws.Range('E33').Value = '  -0.48%  '
$ws.Range('E34').Value = '  +11.25%  '
$ws.Range('E35').Value = '  +5.78%  '
$ws.Range('E36').Value = '  +0.23%  '
$ws.Range('E37').Value = '  +2.81%  '
$ws.Range('E38').Value = '  -0.11%  '
$ws.Range('E39').Value = '  +41.20%  '
$ws.Range('E40').Value = '  -2.78%  '
$ws.Range('E41').Value = '  +5.31%  '
$ws.Range('E42').Value = '  +4.96%  '
$ws.Range('E43').Value = '  +5.65%  '
$ws.Range('E44').Value = '  -0.01%  '
$ws.Range('E45').Value = '  +6.79%  '
$ws.Range('E46').Value = '  +4.34%  '
$ws.Range('E47').Value = '  +15.95%  '
$ws.Range('E48').Value = '  -1.11%  '
$ws.Range('E49').Value = '  +0.35%  '
$ws.Range('E50').Value = '  +2.16%  '
$ws.Range('E51').Value = '  +1.26%  '

# Price column (D) -- force text so Excel does not auto-convert to a number
function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $origStyle = $rng.Style
    $rng.NumberFormat = '@'
    $rng.Value = $val
    $rng.NumberFormat = 'General'
    $rng.Style = $origStyle
}

Set-TextValue 'D2' '96.337.73'
Set-TextValue 'D3' '3.658.28'
Set-TextValue 'D5' '241.57'
Set-TextValue 'D6' '1.89'
Set-TextValue 'D7' '658.19'
Set-TextValue 'D8' '0.424'
Set-TextValue 'D10' '0.999'
Set-TextValue 'D11' '3.658.40'
Set-TextValue 'D12' '44.75'
Set-TextValue 'D14' '6.66'
Set-TextValue 'D15' '4.337.74'
Set-TextValue 'D16' '0.0000270'
Set-TextValue 'D17' '96.047.36'
Set-TextValue 'D18' '8.88'
Set-TextValue 'D19' '3.660.52'
Set-TextValue 'D20' '12.72'
Set-TextValue 'D21' '18.25'
Set-TextValue 'D22' '0.531'
Set-TextValue 'D23' '521.07'
Set-TextValue 'D24' '3.45'
Set-TextValue 'D25' '0.0000204'
Set-TextValue 'D26' '6.88'
Set-TextValue 'D27' '102.18'
Set-TextValue 'D30' '12.24'
Set-TextValue 'D31' '3.02'
Set-TextValue 'D32' '1.00'
Set-TextValue 'D35' '33.20'
Set-TextValue 'D38' '625.64'
Set-TextValue 'D39' '46.43'
Set-TextValue 'D40' '8.72'
Set-TextValue 'D42' '0.958'
Set-TextValue 'D45' '6.24'
Set-TextValue 'D46' '0.0451'
Set-TextValue 'D47' '0.421'
Set-TextValue 'D49' '23.61'
Set-TextValue 'D50' '8.51'
Set-TextValue 'D51' '3.58'
